# Update "想去人数" (F) / "最低票价" (G) figures with freshly scraped counts
# (gh-pages data refresh, commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3830
$ws.Range("G3").Value = 70
$ws.Range("F5").Value = 1368
$ws.Range("G5").Value = 75
$ws.Range("F6").Value = 3828
$ws.Range("F10").Value = 8610
$ws.Range("F11").Value = 485
$ws.Range("F14").Value = 116
$ws.Range("F16").Value = 334
$ws.Range("F17").Value = 91
$ws.Range("F18").Value = 355
$ws.Range("F19").Value = 10894
$ws.Range("F20").Value = 291
$ws.Range("F22").Value = 389
$ws.Range("F28").Value = 2679
$ws.Range("F29").Value = 2077
$ws.Range("F32").Value = 2119
$ws.Range("F33").Value = 894
$ws.Range("F34").Value = 4080
$ws.Range("F35").Value = 2564
$ws.Range("F37").Value = 2582
$ws.Range("F38").Value = 3026
$ws.Range("F39").Value = 1246
$ws.Range("F40").Value = 170
$ws.Range("F41").Value = 755
$ws.Range("F42").Value = 341
$ws.Range("F43").Value = 319
$ws.Range("F45").Value = 107
$ws.Range("F46").Value = 125
$ws.Range("F47").Value = 89

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 15
$ws.Range("F7").Value = 46
$ws.Range("F16").Value = 12
$ws.Range("F22").Value = 42
$ws.Range("F24").Value = 26

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 11

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3830
$ws.Range("G3").Value = 70
$ws.Range("F6").Value = 1368
$ws.Range("G6").Value = 75
$ws.Range("F7").Value = 3828
$ws.Range("F11").Value = 8610
$ws.Range("F12").Value = 485
$ws.Range("F13").Value = 116
$ws.Range("F15").Value = 334
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 355
$ws.Range("F18").Value = 10894
$ws.Range("F19").Value = 291
$ws.Range("F22").Value = 389
$ws.Range("F29").Value = 2679
$ws.Range("F30").Value = 2077
$ws.Range("F31").Value = 2119
$ws.Range("F32").Value = 894
$ws.Range("F34").Value = 4080
$ws.Range("F35").Value = 2564
$ws.Range("F37").Value = 2582
$ws.Range("F38").Value = 3026
$ws.Range("F39").Value = 42
$ws.Range("F40").Value = 1246
$ws.Range("F41").Value = 170
$ws.Range("F42").Value = 755
$ws.Range("F43").Value = 341
$ws.Range("F44").Value = 319
$ws.Range("F45").Value = 107
$ws.Range("F46").Value = 125
$ws.Range("F47").Value = 89
